$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Utilisateurs")

# Add a new user row (16) mirroring the existing rows: Email / FirstName / LastName / State / Password
$ws.Range("B16").Value = "dimitri@fwa.eu"
$ws.Range("C16").Value = "Dimitri"
$ws.Range("D16").Value = "ASHIKHMIN"
$ws.Range("E16").Value = "Active"
$ws.Range("F16").Value = "test"

# Add hyperlink for the new email cell
$ws.Hyperlinks.Add($ws.Range("B16"), "mailto:dimitri@fwa.eu", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "dimitri@fwa.eu")

# Copy formatting from the row above (15) so styles (fills/borders/hyperlink style) match
$ws.Range("A15:F15").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("F16").Select()
